# Rename the worksheet to reflect the "Feedbacks" data now stored in it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Feedbacks of 9 August 2024"

# Update the header row: the meal-tracking columns become feedback columns.
$ws.Range("C1").Value = "Designation"
$ws.Range("D1").Value = "Company"
$ws.Range("E1").Value = "Rating"
$ws.Range("F1").Value = "Description"

# Update the existing data row (row 2) with the new feedback values.
$ws.Range("C2").Value = "HR Manager"
$ws.Range("D2").Value = "Contour Software"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "This is good feedback"

# Append new feedback entries submitted by other users (added on Enter press).
$ws.Range("A3").Value = 28623
$ws.Range("B3").Value = "Sara Akbar"
$ws.Range("C3").Value = "QA Engineer"
$ws.Range("D3").Value = "Contour Software"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "The food was really good today"

$ws.Range("A4").Value = 52837
$ws.Range("B4").Value = "Emma Brown"
$ws.Range("C4").Value = "QA Engineer"
$ws.Range("D4").Value = "Autosoft"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "biryani was good."
